$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string (cell A1) ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 04:22"

# --- Row 38 (Corea del Sur): refresh case counters ---
$ws.Range("B38").Value = 10761
$ws.Range("C38").Value = 9
$ws.Range("D38").Value = 8922
$ws.Range("E38").Value = 1593
$ws.Range("G38").Value = 2
$ws.Range("H38").Value = 246

# --- Rows 130/131: Paraguay overtakes Gabon in the ranking ---
# Row 130 now shows Paraguay with its refreshed totals ...
$ws.Range("A130").Value = "Paraguay"
$ws.Range("B130").Value = 239
$ws.Range("C130").Value = 9
$ws.Range("D130").Value = 102
$ws.Range("E130").Value = 128
$ws.Range("H130").Value = 9

# ... and row 131 now shows Gabon with its previous totals.
$ws.Range("A131").Value = "Gabon"
$ws.Range("B131").Value = 238
$ws.Range("D131").Value = 53
$ws.Range("E131").Value = 182
$ws.Range("H131").Value = 3
